# Append a new paragraph "Klingler" right after the existing "Josh" paragraph.
#
# A plain InsertParagraphAfter()/TypeParagraph() on the existing range causes
# the new paragraph to inherit a <w:pPr><w15:collapsed .../></w:pPr> from the
# paragraph-insertion machinery. The target revision's new paragraph has no
# <w:pPr> at all, so instead we splice in a minimal OOXML fragment for the
# paragraph via Range.InsertXML at the very end of the document's main story,
# which inserts exactly the requested markup without any inherited paragraph
# properties.

$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)

$klinglerParagraphXml = '<?xml version="1.0" standalone="yes"?>' + "`r`n" + `
    '<?mso-application progid="Word.Document"?>' + "`r`n" + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
            '<pkg:xmlData>' + `
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                    '<w:body>' + `
                        '<w:p><w:r><w:t>Klingler</w:t></w:r></w:p>' + `
                    '</w:body>' + `
                '</w:document>' + `
            '</pkg:xmlData>' + `
        '</pkg:part>' + `
    '</pkg:package>'

$end.InsertXML($klinglerParagraphXml)
